$d = $word.ActiveDocument

# The document's headers carry the BTEC logo (currently exported as
# "image1.jpg") and the footers carry the Pearson Edexcel logo (currently
# exported as "image2.png"). The commit renames the embedded picture
# objects' display names: the BTEC logo picture becomes "image2.jpg" and
# the Pearson logo picture becomes "image1.png" (both in every
# header/footer that repeats the logo, i.e. the primary header/footer and
# the "even page" header/footer that mirrors it).

foreach ($sec in $d.Sections) {
    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($k = 1; $k -le $shapes.Count; $k++) {
                $shp = $shapes.Item($k)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
    for ($f = 1; $f -le 3; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($k = 1; $k -le $shapes.Count; $k++) {
                $shp = $shapes.Item($k)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
